$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("D6").Value = "2016-21-19 22:21:05"
$ws.Range("B9").Value = "Ready for handoff"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("D9").Value = "2016-21-19 22:21:05"
$ws.Range("D10").Value = "2016-21-19 22:21:05"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E6").Value = "2016-03-19 22:21:01"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("E9").Value = "2016-03-19 22:21:01"
$ws.Range("E10").Value = "2016-03-19 22:21:01"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("E6").Value = "2016-03-19 22:21:05"
$ws.Range("C9").Value = "Ready for handoff"
$ws.Range("E9").Value = "2016-03-19 22:21:05"
$ws.Range("E10").Value = "2016-03-19 22:21:05"
